$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 99.42857
$ws.Range("I2").Value = 99.2
$ws.Range("K2").Value = 99.2
$ws.Range("M2").Value = 13.8
$ws.Range("H18").Value = 2099.25
$ws.Range("I18").Value = 3149.5
$ws.Range("J18").Value = 1049.0
$ws.Range("K18").Value = 3149.5
$ws.Range("L18").Value = 1049.0
$ws.Range("M18").Value = -2865.5
$ws.Range("N18").Value = -1617.0
$ws.Range("H55").Value = 145.66667
$ws.Range("I55").Value = 145.66667
$ws.Range("K55").Value = 145.66667
$ws.Range("M55").Value = 68.33332999999999
$ws.Range("H58").Value = 1472.2858
$ws.Range("I58").Value = 51.0
$ws.Range("J58").Value = 10000.0
$ws.Range("K58").Value = 153.0
$ws.Range("L58").Value = 30000.0
$ws.Range("M58").Value = -3.0
$ws.Range("N58").Value = -30300.0
$ws.Range("H97").Value = 2084.8333
$ws.Range("J97").Value = 921.8
$ws.Range("L97").Value = 2765.4
$ws.Range("N97").Value = -3757.4
$ws.Range("H98").Value = 2220.4546
$ws.Range("J98").Value = 2228.6667
$ws.Range("L98").Value = 2228.6667
$ws.Range("N98").Value = -5224.6667
$ws.Range("H112").Value = 4200.0
$ws.Range("I112").Value = 4000.0
$ws.Range("J112").Value = 4250.0
$ws.Range("K112").Value = 12000.0
$ws.Range("L112").Value = 12750.0
$ws.Range("M112").Value = -10892.0
$ws.Range("N112").Value = -14966.0
$ws.Range("H113").Value = 4334.1665
$ws.Range("I113").Value = 4001.25
$ws.Range("J113").Value = 5000.0
$ws.Range("K113").Value = 4001.25
$ws.Range("L113").Value = 5000.0
$ws.Range("M113").Value = -747.25
$ws.Range("N113").Value = -11508.0
$ws.Range("H116").Value = 41093.43
$ws.Range("I116").Value = 9284.667
$ws.Range("J116").Value = 64950.0
$ws.Range("K116").Value = 9284.667
$ws.Range("L116").Value = 64950.0
$ws.Range("M116").Value = -5842.666999999999
$ws.Range("N116").Value = -71834.0
$ws.Range("H122").Value = 2220.4546
$ws.Range("J122").Value = 2228.6667
$ws.Range("L122").Value = 6686.000100000001
$ws.Range("N122").Value = -11586.0001
$ws.Range("H127").Value = 3196.2727
$ws.Range("I127").Value = 3219.5173
$ws.Range("J127").Value = 3027.75
$ws.Range("K127").Value = 9658.5519
$ws.Range("L127").Value = 9083.25
$ws.Range("M127").Value = -4698.5519
$ws.Range("N127").Value = -19003.25
$ws.Range("H129").Value = 1824.5769
$ws.Range("I129").Value = 1645.2632
$ws.Range("J129").Value = 2311.2856
$ws.Range("K129").Value = 4935.7896
$ws.Range("L129").Value = 6933.8568
$ws.Range("M129").Value = 64.21039999999994
$ws.Range("N129").Value = -16933.8568
$ws.Range("H136").Value = 198858.67
$ws.Range("J136").Value = 198858.67
$ws.Range("L136").Value = 198858.67
$ws.Range("N136").Value = -209058.67
$ws.Range("H137").Value = 22227570.0
$ws.Range("I137").Value = 66667892.0
$ws.Range("J137").Value = 7409.8
$ws.Range("K137").Value = 200003676.0
$ws.Range("L137").Value = 22229.4
$ws.Range("M137").Value = -200001126.0
$ws.Range("N137").Value = -27329.4
$ws.Range("H138").Value = 2704.2327
$ws.Range("J138").Value = 3604.3333
$ws.Range("L138").Value = 10812.9999
$ws.Range("N138").Value = -21092.9999
$ws.Range("H139").Value = 77863.63
$ws.Range("J139").Value = 77863.63
$ws.Range("L139").Value = 77863.63
$ws.Range("N139").Value = -88143.63
$ws.Range("H140").Value = 100000.0
$ws.Range("J140").Value = 100000.0
$ws.Range("L140").Value = 100000.0
$ws.Range("N140").Value = -110360.0
$ws.Range("H141").Value = 4907.6665
$ws.Range("J141").Value = 14142.857
$ws.Range("L141").Value = 42428.571
$ws.Range("N141").Value = -52788.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 181864.05
$ws.Range("I32").Value = 254402.88
$ws.Range("K32").Value = 254402.88
$ws.Range("M32").Value = -254115.88
$ws.Range("H61").Value = 2328911.8
$ws.Range("I61").Value = 3254.4324
$ws.Range("J61").Value = 16670465.0
$ws.Range("K61").Value = 3254.4324
$ws.Range("L61").Value = 16670465.0
$ws.Range("M61").Value = -3042.4324
$ws.Range("N61").Value = -16670889.0
$ws.Range("H122").Value = 1979.1818
$ws.Range("I122").Value = 1598.0
$ws.Range("K122").Value = 4794.0
$ws.Range("M122").Value = -2344.0
$ws.Range("H136").Value = 2328911.8
$ws.Range("I136").Value = 3254.4324
$ws.Range("J136").Value = 16670465.0
$ws.Range("K136").Value = 9763.2972
$ws.Range("L136").Value = 50011395.0
$ws.Range("M136").Value = -7213.297200000001
$ws.Range("N136").Value = -50016495.0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 95.0
$ws.Range("I22").Value = 95.0
$ws.Range("K22").Value = 95.0
$ws.Range("M22").Value = 78.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4633353.5
$ws.Range("I31").Value = 4633353.5
$ws.Range("K31").Value = 4633353.5
$ws.Range("M31").Value = -4633058.5
$ws.Range("H34").Value = 4633353.5
$ws.Range("I34").Value = 4633353.5
$ws.Range("K34").Value = 4633353.5
$ws.Range("M34").Value = -4633151.5
$ws.Range("H64").Value = 0.0
$ws.Range("J64").Value = 0.0
$ws.Range("L64").Value = 0.0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0.0
$ws.Range("J67").Value = 0.0
$ws.Range("L67").Value = 0.0
$ws.Range("N67").ClearContents()
$ws.Range("H134").Value = 1856.119
$ws.Range("I134").Value = 1491.7567
$ws.Range("K134").Value = 4475.2701
$ws.Range("M134").Value = -1940.2701
$ws.Range("H140").Value = 86373.81
$ws.Range("J140").Value = 86373.81
$ws.Range("L140").Value = 86373.81
$ws.Range("N140").Value = -96733.81

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 998.0
$ws.Range("I22").Value = 998.0
$ws.Range("J22").Value = 0.0
$ws.Range("K22").Value = 2994.0
$ws.Range("L22").Value = 0.0
$ws.Range("M22").Value = -2825.0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 998.0
$ws.Range("I27").Value = 998.0
$ws.Range("J27").Value = 0.0
$ws.Range("K27").Value = 2994.0
$ws.Range("L27").Value = 0.0
$ws.Range("M27").Value = -2892.0
$ws.Range("N27").ClearContents()
$ws.Range("H107").Value = 1685.7812
$ws.Range("I107").Value = 371.0
$ws.Range("K107").Value = 1113.0
$ws.Range("M107").Value = 807.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 820.0476
$ws.Range("I16").Value = 792.4375
$ws.Range("K16").Value = 792.4375
$ws.Range("M16").Value = -622.4375
$ws.Range("H46").Value = 3689.8667
$ws.Range("J46").Value = 5660.8887
$ws.Range("L46").Value = 5660.8887
$ws.Range("N46").Value = -6036.8887
$ws.Range("H82").Value = 757.6667
$ws.Range("I82").Value = 788.5714
$ws.Range("K82").Value = 788.5714
$ws.Range("M82").Value = -427.5714
$ws.Range("H85").Value = 757.6667
$ws.Range("I85").Value = 788.5714
$ws.Range("K85").Value = 788.5714
$ws.Range("M85").Value = 459.4286
$ws.Range("H136").Value = 8935527.0
$ws.Range("I136").Value = 6253304.5
$ws.Range("J136").Value = 15641082.0
$ws.Range("K136").Value = 18759913.5
$ws.Range("L136").Value = 46923246.0
$ws.Range("M136").Value = -18757363.5
$ws.Range("N136").Value = -46928346.0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 607141.9
$ws.Range("I4").Value = 607141.9
$ws.Range("K4").Value = 607141.9
$ws.Range("M4").Value = -607028.9
$ws.Range("H32").Value = 30000.0
$ws.Range("I32").Value = 30000.0
$ws.Range("K32").Value = 30000.0
$ws.Range("M32").Value = -29683.0
$ws.Range("H34").Value = 0.0
$ws.Range("J34").Value = 0.0
$ws.Range("L34").Value = 0.0
$ws.Range("N34").ClearContents()
$ws.Range("H38").Value = 500.0
$ws.Range("J38").Value = 500.0
$ws.Range("L38").Value = 500.0
$ws.Range("N38").Value = -1446.0
$ws.Range("H44").Value = 0.0
$ws.Range("I44").Value = 0.0
$ws.Range("K44").Value = 0.0
$ws.Range("M44").ClearContents()
$ws.Range("H49").Value = 25250000.0
$ws.Range("J49").Value = 0.0
$ws.Range("L49").Value = 0.0
$ws.Range("N49").ClearContents()
$ws.Range("H58").Value = 4994.6665
$ws.Range("I58").Value = 4994.6665
$ws.Range("K58").Value = 4994.6665
$ws.Range("M58").Value = -4686.6665
$ws.Range("H122").Value = 80145.47
$ws.Range("I122").Value = 7235.222
$ws.Range("K122").Value = 21705.666
$ws.Range("M122").Value = -19255.666
$ws.Range("H126").Value = 4047.4546
$ws.Range("I126").Value = 4778.4287
$ws.Range("K126").Value = 14335.2861
$ws.Range("M126").Value = -11865.2861
$ws.Range("H132").Value = 27780512.0
$ws.Range("I132").Value = 33335754.0
$ws.Range("K132").Value = 100007262.0
$ws.Range("M132").Value = -100004732.0
